$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Claro"
$ws.Range("B2").Value = "250,00"
$ws.Range("C2").Value = "'12/06/2023"
$ws.Range("D2").Value = "Um exemplo!`n"
